$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''62.045.14'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '''  -0.09%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = '''3.415.14'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '''  -0.05%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = '''0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '''  -0.14%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = '''408.76'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Value = '''129.35'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '''  -3.53%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = '''0.639'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '''  +7.61%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = '''  -0.06%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = '''0.733'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '''  +7.27%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('E10').Value = '''  +19.32%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = '''42.54'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '''  -0.59%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = '''0.0000222'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '''  +70.90%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('E13').Value = '''  -0.44%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = '''3.958.00'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '''  -0.05%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = '''8.91'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '''  +5.60%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = '''20.76'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '''  +4.41%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = '''3.421.13'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '''  -0.31%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = '''12.14'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '''  +10.13%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('E19').Value = '''  +5.66%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = '''61.915.95'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '''  -0.36%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = '''410.92'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '''  +31.16%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = '''89.46'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '''  +6.24%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('E23').Value = '''  -0.43%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = '''13.07'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '''  +1.48%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = '''3.24'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '''  +2.05%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = '''33.07'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '''  +11.59%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = '''8.83'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '''  +7.89%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = '''4.77'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '''  -0.42%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = '''7.59'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '''  -0.28%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('E30').Value = '''  -1.64%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('B31').Value = '''Hedera'
$ws.Range('B31').Style = 'Normal'
$ws.Range('C31').Value = '''https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('C31').Style = 'Normal'
$ws.Range('D31').Value = '''0.116'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '''  +0.36%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('B32').Value = '''Cosmos'
$ws.Range('B32').Style = 'Normal'
$ws.Range('C32').Value = '''https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('C32').Style = 'Normal'
$ws.Range('D32').Value = '''11.87'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '''  +4.12%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('B33').Value = '''Kaspa'
$ws.Range('B33').Style = 'Normal'
$ws.Range('C33').Value = '''https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('C33').Style = 'Normal'
$ws.Range('D33').Value = '''0.171'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '''  -2.04%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = '''42.77'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '''  -0.29%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = '''  +0.80%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').Value = '''  +3.36%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = '''54.27'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '''  +4.57%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = '''0.997'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '''  -0.22%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('B39').Value = '''Stellar'
$ws.Range('B39').Style = 'Normal'
$ws.Range('C39').Value = '''https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('C39').Style = 'Normal'
$ws.Range('D39').Value = '''0.135'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '''  +7.56%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('B40').Value = '''LidoDAOToken'
$ws.Range('B40').Style = 'Normal'
$ws.Range('C40').Value = '''https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('C40').Style = 'Normal'
$ws.Range('D40').Value = '''3.35'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '''  -2.10%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = '''2.92'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '''  -1.38%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = '''0.311'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '''  +4.19%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = '''141.65'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '''  +3.35%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('E44').Value = '''  -1.15%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = '''4.12'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '''  +1.83%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('E46').Value = '''  +8.69%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = '''16.65'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '''  -0.77%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = '''21.96'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '''  +2.24%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = '''2.112.90'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '''  -0.45%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = '''  +3.14%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('E51').Value = '''  +15.31%  '
$ws.Range('E51').Style = 'Normal'
